# Ver-3.0.2 Centrality metrics of multi-graphs have been changed to `unweighted`.
# Update the Degree-centrality value columns (O, Q, S) on rows 2-4 from raw
# (weighted) degree counts to normalized (unweighted) degree centrality values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("O2").Value = 0.2779661016949153
$ws.Range("Q2").Value = 0.2711864406779661
$ws.Range("S2").Value = 0.2372881355932203

$ws.Range("O3").Value = 0.4414414414414414
$ws.Range("Q3").Value = 0.3963963963963964
$ws.Range("S3").Value = 0.3693693693693694

$ws.Range("O4").Value = 0.5017793594306049
$ws.Range("Q4").Value = 0.4733096085409252
$ws.Range("S4").Value = 0.4270462633451957
